# Populate the rule-number column (A3:A8) on the "WeightClassification"
# DMN decision-table sheet with the rule index (1-6) for each row.
# These cells were previously blank/empty numeric cells and must become
# text values "1".."6" (matching the DMN rule-numbering convention used
# by the regenerated XLSX decision tables).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeightClassification")

# Force the cells to be stored as text (not auto-coerced to numbers) by
# pre-formatting the range as Text before assigning the numeric-looking
# string values.
$ws.Range("A3:A8").NumberFormat = "@"

$ws.Range("A3").Value = "1"
$ws.Range("A4").Value = "2"
$ws.Range("A5").Value = "3"
$ws.Range("A6").Value = "4"
$ws.Range("A7").Value = "5"
$ws.Range("A8").Value = "6"

Write-Host "Updated A3:A8 on WeightClassification with rule numbers 1-6"
